$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its text formatting so values like
# "244.08" or "1.002" are not reinterpreted as numbers/dates.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.833.56'
$ws.Range("E2").Value = '  +1.49%  '

$ws.Range("D3").Value = '1.878.56'
$ws.Range("E3").Value = '  +2.24%  '

$ws.Range("E4").Value = '  +0.51%  '

$ws.Range("D5").Value = '244.08'
$ws.Range("E5").Value = '  -1.29%  '

$ws.Range("E6").Value = '  +0.52%  '

$ws.Range("D7").Value = '0.4943'
$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").Value = '44.34'
$ws.Range("E8").Value = '  +0.20%  '

$ws.Range("D9").Value = '0.2903'
$ws.Range("E9").Value = '  +3.51%  '

$ws.Range("D10").Value = '0.06595'
$ws.Range("E10").Value = '  +2.71%  '

$ws.Range("D11").Value = '1.888.93'
$ws.Range("E11").Value = '  +3.38%  '

$ws.Range("D12").Value = '16.88'
$ws.Range("E12").Value = '  +0.37%  '

$ws.Range("D13").Value = '0.07192'
$ws.Range("E13").Value = '  +1.28%  '

$ws.Range("D14").Value = '0.6688'
$ws.Range("E14").Value = '  +2.74%  '

$ws.Range("D15").Value = '85.43'
$ws.Range("E15").Value = '  +1.11%  '

$ws.Range("D16").Value = '4.813'
$ws.Range("E16").Value = '  +1.79%  '

$ws.Range("D17").Value = '29.855.35'
$ws.Range("E17").Value = '  +1.60%  '

$ws.Range("D18").Value = '0.000007807'
$ws.Range("E18").Value = '  +6.20%  '

$ws.Range("D19").Value = '1.002'
$ws.Range("E19").Value = '  +0.34%  '

$ws.Range("D20").Value = '12.76'
$ws.Range("E20").Value = '  +2.62%  '

$ws.Range("D21").Value = '2.133.51'
$ws.Range("E21").Value = '  +3.51%  '

$ws.Range("E22").Value = '  +0.52%  '

$ws.Range("D23").Value = '4.740'
$ws.Range("E23").Value = '  +3.29%  '

$ws.Range("D24").Value = '5.580'
$ws.Range("E24").Value = '  +2.62%  '

$ws.Range("D25").Value = '9.105'
$ws.Range("E25").Value = '  +2.42%  '

$ws.Range("D26").Value = '147.95'
$ws.Range("E26").Value = '  +2.88%  '

$ws.Range("D27").Value = '134.37'
$ws.Range("E27").Value = '  +1.88%  '

$ws.Range("D28").Value = '16.69'
$ws.Range("E28").Value = '  +1.17%  '

$ws.Range("D29").Value = '1.923'
$ws.Range("E29").Value = '  +0.73%  '

$ws.Range("E30").Value = '  -1.65%  '

$ws.Range("D31").Value = '4.170'
$ws.Range("E31").Value = '  +0.07%  '

$ws.Range("D32").Value = '0.08682'
$ws.Range("E32").Value = '  +3.62%  '

$ws.Range("D33").Value = '3.932'
$ws.Range("E33").Value = '  +3.32%  '

$ws.Range("D34").Value = '0.05065'
$ws.Range("E34").Value = '  +2.23%  '

$ws.Range("D35").Value = '1.108'
$ws.Range("E35").Value = '  +0.19%  '

$ws.Range("D36").Value = '0.7023'
$ws.Range("E36").Value = '  +3.87%  '

$ws.Range("D37").Value = '2.682'
$ws.Range("E37").Value = '  -0.13%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '2.205'
$ws.Range("E38").Value = '  -3.44%  '

$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.688'
$ws.Range("E39").Value = '  -0.99%  '

$ws.Range("D40").Value = '0.9328'
$ws.Range("E40").Value = '  -2.15%  '

$ws.Range("D41").Value = '0.01641'
$ws.Range("E41").Value = '  +3.23%  '

$ws.Range("D42").Value = '6.049'
$ws.Range("E42").Value = '  -3.03%  '

$ws.Range("D43").Value = '0.9998'
$ws.Range("E43").Value = '  +0.13%  '

$ws.Range("D44").Value = '102.67'
$ws.Range("E44").Value = '  +0.03%  '

$ws.Range("D45").Value = '0.4168'
$ws.Range("E45").Value = '  +1.84%  '

$ws.Range("D46").Value = '7.445'
$ws.Range("E46").Value = '  +2.57%  '

$ws.Range("D47").Value = '0.1257'
$ws.Range("E47").Value = '  +2.69%  '

$ws.Range("D48").Value = '0.05714'
$ws.Range("E48").Value = '  +2.54%  '

$ws.Range("D49").Value = '32.51'
$ws.Range("E49").Value = '  +1.78%  '

$ws.Range("D50").Value = '8.183'
$ws.Range("E50").Value = '  +1.05%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '55.81'
$ws.Range("E51").Value = '  +3.57%  '
